# Rename the metadata tab from "column_specs" to "column_names"
# (column specs are now internal, so the tab now only lists raw/new
# column name pairs without the col_type column).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("column_specs")
$ws.Name = "column_names"

# Drop the now-unneeded col_type column from the Table7 table (and its
# backing worksheet range) on the renamed tab.
$lo = $ws.ListObjects.Item("Table7")
$lo.ListColumns.Item("col_type").Delete()
